$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new rows of data to the StatusData sheet (rows 4 and 5).
# Phone numbers are entered with a leading apostrophe so Excel stores
# them as text (matching the existing text-typed mobile numbers in
# rows 2-3) instead of auto-converting them to numeric values.

$ws.Range("A4").Value = "Anita"
$ws.Range("B4").Value = "'8368547181"
$ws.Range("C4").Value = "Trusted"

$ws.Range("A5").Value = "Ayush"
$ws.Range("B5").Value = "'8882292108"
$ws.Range("C5").Value = "Blacklisted"
